# Fruta / hortaliza, semanal
# Re-applies the weekly refresh: the data rows (2-13) are shuffled around
# (same underlying records, new row positions) with the corresponding
# cell values updated in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: row number -> hashtable of column letter -> new value
$updates = @{
    2  = @{ D = 44159; L = "Segunda"; M = 200; N = 6500;                 P = 6750;                                               S = 4500 }
    3  = @{ D = 44523; L = "Primera"; M = 300; N = 3700; O = 3800; P = 3750; Q = "`$/kilo"; R = "Región del Maule";               S = 3750; T = 1 }
    4  = @{ D = 44169;                         N = 5500; O = 6000; P = 5750;                R = "Provincia de Curicó";            S = 3833 }
    5  = @{ D = 44530;            M = 160;     N = 3600; O = 3700; P = 3650;                                                     S = 3650 }
    6  = @{ D = 44537;            M = 400;     N = 5000; O = 5500; P = 5250; Q = "`$/bandeja 12 canastillos 125 gramos";          S = 3500; T = 1.5 }
    7  = @{ D = 44162;            M = 100;     N = 7000; O = 7000; P = 7000;                                                     S = 4667 }
    8  = @{ D = 44162; L = "Segunda"; M = 100; N = 6500; O = 6500; P = 6500; Q = "`$/bandeja 12 canastillos 125 gramos"; R = "Provincia de Curicó"; S = 4333; T = 1.5 }
    9  = @{ D = 44166;            M = 200;     N = 6000; O = 6500; P = 6250;                                                     S = 4167 }
    10 = @{ D = 44516;            M = 80;      N = 3700; O = 3800; P = 3750;                                                     S = 3750 }
    11 = @{ D = 44533; L = "Primera"; M = 400; N = 3500; O = 3600; P = 3550; Q = "`$/kilo"; R = "Región del Maule";               S = 3550; T = 1 }
    12 = @{ D = 44176;            M = 300;     N = 5000; O = 6000; P = 5500; Q = "`$/bandeja 12 canastillos 125 gramos"; R = "Provincia de Curicó"; S = 3667; T = 1.5 }
    13 = @{ D = 44519;                         N = 3700; O = 3800; P = 3750; Q = "`$/kilo"; R = "Región del Maule";               S = 3750; T = 1 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
